# Applies the Quiz.xlsx edit: module 1 (B1..B10) questions replaced by
# module 2 (C1..C10) questions, and the "correct" answer letters updated
# for the rows whose index changed in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (module number): keep cell type as Text ("1" -> "2") ---
# Setting NumberFormat to Text before assigning the value keeps the cell
# a shared-string "2" instead of Excel auto-converting it to a Number;
# resetting Style afterwards drops the now-unneeded per-cell style so the
# cell ends up indistinguishable (style-wise) from the original.
$moduleRange = $ws.Range("A2:A11")
$moduleRange.NumberFormat = "@"
$moduleRange.Value = '2'
$moduleRange.Style = "Normal"

# --- Row 2: question C1 ---
$ws.Range("B2").Value = 'C1'  # id
$ws.Range("C2").Value = 'Estas a punto de entrar y vas tarde. Tenes 45 segundos. Cual preparacion minima te deja mejor parado para no improvisar adentro?'  # question
$ws.Range("D2").Value = 'Entrar directo y preguntar “¿qué te anoto?” para ahorrar tiempo, sin preparar nada afuera.'  # a
$ws.Range("E2").Value = 'Mirar solo promos generales y entrar a ofrecerlas, para no meter presión ni discutir.'  # b
$ws.Range("F2").Value = 'Marcar GPS, mirar en Bot/App 2 datos (compra fija + última compra) y entrar con plan simple: base + 1 idea.'  # c
$ws.Range("G2").Value = 'Pensar una estrategia sin mirar Bot/App: “adentro veo” y arrancar. Y seguís con el pedido en el sistema.'  # d
$ws.Range("H2").Value = 'Llamar al cliente antes de entrar para preguntarle qué quiere y armar el pedido por teléfono.'  # e
$ws.Range("I2").Value = 'c'  # correct

# --- Row 3: question C2 ---
$ws.Range("B3").Value = 'C2'  # id
$ws.Range("C3").Value = 'En el Bot ves que un cliente que siempre compraba un producto fuerte bajo mucho el volumen las ultimas 2 visitas. Que preparación es más inteligente antes de entrar?'  # question
$ws.Range("D3").Value = 'Ignorarlo y sostener la visita igual, asumiendo que ya no le interesa ese producto.'  # a
$ws.Range("E3").Value = 'Entrar ofreciendo solo productos nuevos para subir ticket, sin revisar qué pasó con el fuerte.'  # b
$ws.Range("F3").Value = 'Preparar una pregunta corta para entender el motivo (rotación/precio/competencia) y ajustar reposición a su realidad.'  # c
$ws.Range("G3").Value = 'Entrar a reclamarle por qué bajó y presionarlo para que vuelva a subir el volumen.'  # d
$ws.Range("H3").Value = 'No visitarlo por ahora, asumir que no vale la pena y pasar al siguiente cliente hoy.'  # e
$ws.Range("I3").Value = 'c'  # correct

# --- Row 4: question C3 ---
$ws.Range("B4").Value = 'C3'  # id
$ws.Range("C4").Value = 'Cual de estas cosas NO es parte del checklist de preparacion inmediata antes de entrar?'  # question
$ws.Range("D4").Value = 'Marcar GPS al llegar. Para no meter presión ni discutir.'  # a
$ws.Range("E4").Value = 'Revisar qué compra siempre, cantidades y frecuencia. Para que quede corto y claro.'  # b
$ws.Range("F4").Value = 'Identificar si compra categorías del portfolio a otros proveedores.'  # c
$ws.Range("G4").Value = 'Memorizar una lista completa de 10 promos para decirlas todas'  # d
$ws.Range("H4").Value = 'Recordar su problema principal de la visita anterior.'  # e
$ws.Range("I4").Value = 'd'  # correct

# --- Row 5: question C4 ---
$ws.Range("B5").Value = 'C4'  # id
$ws.Range("C5").Value = 'Cliente nuevo sin historial. Cual seria una preparacion minima aceptable antes de entrar?'  # question
$ws.Range("D5").Value = 'Entrar sin preparar nada y empezar la conversación “a ver qué sale” en el momento.'  # a
$ws.Range("E5").Value = 'Marcar GPS, revisar ofertas zonales y armar un pedido base típico más una pregunta de diagnóstico.'  # b
$ws.Range("F5").Value = 'Marcar GPS y confiar en improvisar todo adentro, sin revisar información previa del recorrido.'  # c
$ws.Range("G5").Value = 'Entrar con un discurso de promos generales y esperar que alguna alternativa le interese hoy.'  # d
$ws.Range("H5").Value = 'Postergar la visita hasta tener compras registradas y recién ahí empezar a trabajarlo en serio.'  # e
$ws.Range("I5").Value = 'b'  # correct

# --- Row 6: question C5 ---
$ws.Range("B6").Value = 'C5'  # id
$ws.Range("C6").Value = 'Tenes 6 promos activas. En el Bot ves que el cliente viene quejandose de margen. Que decision previa es mas solida?'  # question
$ws.Range("D6").Value = 'Entrar y recitar las 6 promos seguidas, esperando que alguna le cierre en el momento.'  # a
$ws.Range("E6").Value = 'Elegir la promo más fuerte del mes, aunque no conecte con el problema de margen del cliente.'  # b
$ws.Range("F6").Value = 'Evitar promos por completo y limitarte a tomar el pedido habitual sin proponer alternativas.'  # c
$ws.Range("G6").Value = 'Seleccionar 1–2 promos que mejoren margen y preparar una frase simple de beneficio para ese cliente.'  # d
$ws.Range("H6").Value = 'Ofrecer solo lanzamientos nuevos y dejar afuera promos pensadas para margen y rotación del local.'  # e
$ws.Range("I6").Value = 'd'  # correct

# --- Row 7: question C6 ---
$ws.Range("B7").Value = 'C6'  # id
$ws.Range("C7").Value = 'Antes de entrar, recordas que el cliente te dijo la visita pasada: ''no tengo mas lugar''. Que mini estrategia pre-armas?'  # question
$ws.Range("D7").Value = 'Entrar con foco en subir volumen para aprovechar el viaje, aunque ya avisó que no tiene espacio.'  # a
$ws.Range("E7").Value = 'Entrar con foco en liberar espacio: ver qué no rota y proponer alta rotación en pocas unidades y formato chico.'  # b
$ws.Range("F7").Value = 'No ofrecer nada nuevo y pedir solo “lo de siempre”, aunque la traba de espacio siga igual.'  # c
$ws.Range("G7").Value = 'Cambiar el tema hacia impuestos y desviar la conversación, sin resolver la traba principal del cliente.'  # d
$ws.Range("H7").Value = 'Evitar mirar depósito/estantería para no incomodar y seguir la visita como si no hubiera problema.'  # e
$ws.Range("I7").Value = 'b'  # correct

# --- Row 8: question C7 ---
$ws.Range("B8").Value = 'C7'  # id
$ws.Range("C8").Value = 'Te diste cuenta de que entraste sin marcar GPS. Que accion es mas profesional?'  # question
$ws.Range("D8").Value = 'Marcar GPS apenas puedas (si hace falta, salir un segundo a la puerta) para que la visita quede registrada.'  # a
$ws.Range("E8").Value = 'Seguir la visita normal sin marcar GPS y darlo por perdido en esta vuelta.'  # b
$ws.Range("F8").Value = 'Marcar el GPS al final del día desde cualquier lugar, solo para que figure en el sistema.'  # c
$ws.Range("G8").Value = 'Pedirle al cliente que confirme que estuviste y usar eso como reemplazo del GPS.'  # d
$ws.Range("H8").Value = 'Marcar GPS en el negocio siguiente “para compensar” y cerrar el número de visitas.'  # e
$ws.Range("I8").Value = 'a'  # correct

# --- Row 9: question C8 ---
$ws.Range("B9").Value = 'C8'  # id
$ws.Range("C9").Value = 'Ordena la secuencia mas logica antes de entrar (1 minuto de preparacion).'  # question
$ws.Range("D9").Value = 'Marcar GPS -> mirar Bot/App -> pensar 1–2 ideas -> entrar y saludar con foco.'  # a
$ws.Range("E9").Value = 'Entrar y saludar -> marcar GPS -> mirar Bot/App -> pensar estrategia ya adentro.'  # b
$ws.Range("F9").Value = 'Mirar Bot/App -> entrar y saludar -> marcar GPS -> pensar estrategia a mitad de visita.'  # c
$ws.Range("G9").Value = 'Pensar estrategia -> entrar y saludar -> mirar Bot/App -> marcar GPS al final.'  # d
$ws.Range("H9").Value = 'Marcar GPS -> entrar -> pensar estrategia sin mirar Bot/App y salir improvisando.'  # e
$ws.Range("I9").Value = 'a'  # correct

# --- Row 10: question C9 ---
$ws.Range("B10").Value = 'C9'  # id
$ws.Range("C10").Value = 'Cual de estas senales en el Bot/App deberia prenderte una luz para preparar la visita con mas cuidado?'  # question
$ws.Range("D10").Value = 'Compra siempre lo mismo con la misma frecuencia, y todo viene estable, sin cambios grandes.'  # a
$ws.Range("E10").Value = 'Tiene pedidos regulares y un comportamiento parejo, sin cambios relevantes. Para seguir con la visita con orden.'  # b
$ws.Range("F10").Value = 'Tiene histórico corto pero constante, sin variaciones fuertes, y mantiene el mismo patrón.'  # c
$ws.Range("G10").Value = 'No tiene promos activas hoy, pero su compra viene normal y no hay alertas de caída.'  # d
$ws.Range("H10").Value = 'Bajó de golpe un producto que antes rotaba bien, y se repitió en las últimas dos visitas.'  # e
$ws.Range("I10").Value = 'e'  # correct

# --- Row 11: question C10 ---
$ws.Range("B11").Value = 'C10'  # id
$ws.Range("C11").Value = 'Cual de estas notas de pre-visita es la mas util (1 renglon) para entrar con plan?'  # question
$ws.Range("D11").Value = '“Repasar qué rotó la última visita y entrar con 1 pregunta corta para entender qué cambió.”'  # a
$ws.Range("E11").Value = '“Priorizar 2 promos que encajen con su compra habitual y mencionarlas si aparece la oportunidad.”'  # b
$ws.Range("F11").Value = '“Repasar caballitos X e Y y chequear si hay quiebre/stock para ajustar reposición sin inflar.”'  # c
$ws.Range("G11").Value = '“Entrar a escuchar primero y, según lo que diga, armar una propuesta simple en el momento.”'  # d
$ws.Range("H11").Value = '“Repongo X e Y (caballitos) + propongo prueba acotada para margen o ahorro de viajes hoy.” Respuestas + explicacion (por que es correcta)'  # e
$ws.Range("I11").Value = 'e'  # correct

